$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column test results for Microdol 1 (E) and Microdol 5 (F)
$ws.Range("E6").Value = 66.099999999999994
$ws.Range("F6").Value = 81.3

$ws.Range("E7").Value = 45.2
$ws.Range("F7").Value = 50.1

$ws.Range("E8").Value = 32.799999999999997
$ws.Range("F8").Value = 40.1

$ws.Range("E9").Value = 29.3
$ws.Range("F9").Value = 33.6

$ws.Range("E10").Value = 29.6
$ws.Range("F10").Value = 35.6

# Update selection on the sheet view to E24
$ws.Range("E24").Select()
